# Atualizei os dados da bibi e da add
# - Corrige o valor de faturamento do dia 25 (05/2025), linha 26
# - Insere um novo registro (dia 26 de 05/2025) logo apos a linha 26,
#   empurrando as linhas seguintes uma posicao para baixo

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrige o total_venda da linha 26 (dia 25, mes 05/2025)
$ws.Range("B26").Value = 3824.9

# Insere uma nova linha logo apos a linha 26, deslocando as demais para baixo
$ws.Rows.Item(27).Insert()

# Preenche os dados do novo registro inserido (dia 26, 05/2025)
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 27393.77
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 2025
$ws.Range("E27").Value = "05/2025"
